# Enabling all testcases ECTEST
# Set the Runmode column (E) to "Yes" for every testcase row (2-33),
# enabling all test cases. Then update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value -ne "Yes") {
        $cell.Value = "Yes"
    }
}

$ws.Range("E2:E33").Select() | Out-Null
